# Add 2022-Q4 data: insert a new worksheet with the quarterly fund holdings
# table, and add the corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: Insert a new worksheet "2022-Q4" right after "总计" (i.e. right
# before the existing "2022-Q3" sheet), matching the target tab order:
#   总计, 2022-Q4, 2022-Q3, 2022-Q1, 2021-Q2, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------
$existing = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Add($existing)
$ws.Name = "2022-Q4"

function Set-HeaderCell($cell, $text) {
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Set-IndexCell($cell, $n) {
    $cell.Value = $n
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Header row (row 1), columns B..H
Set-HeaderCell $ws.Cells.Item(1, 2) "基金代码"
Set-HeaderCell $ws.Cells.Item(1, 3) "基金名称"
Set-HeaderCell $ws.Cells.Item(1, 4) "基金规模"
Set-HeaderCell $ws.Cells.Item(1, 5) "股票总仓位"
Set-HeaderCell $ws.Cells.Item(1, 6) "仓位占比"
Set-HeaderCell $ws.Cells.Item(1, 7) "持有市值(亿元)"
Set-HeaderCell $ws.Cells.Item(1, 8) "仓位排名"

# Data rows 2..5 : index, code, name, size, total position, position pct,
# held market value, position rank
$rows = @(
    @(0, "007553", "中信建投医改灵活配置混合C", "11.40", "95.02", "3.76", "0.4286", 10),
    @(1, "002408", "中信建投医改灵活配置混合A", "11.28", "95.02", "3.76", "0.4241", 10),
    @(2, "010090", "中信建投医药健康混合A",     "2.95",  "95.01", "3.69", "0.1089", 9),
    @(3, "010091", "中信建投医药健康混合C",     "1.98",  "95.01", "3.69", "0.0731", 9)
)

$r = 2
foreach ($row in $rows) {
    Set-IndexCell $ws.Cells.Item($r, 1) $row[0]
    Set-TextCell  $ws.Cells.Item($r, 2) $row[1]
    Set-TextCell  $ws.Cells.Item($r, 3) $row[2]
    Set-TextCell  $ws.Cells.Item($r, 4) $row[3]
    Set-TextCell  $ws.Cells.Item($r, 5) $row[4]
    Set-TextCell  $ws.Cells.Item($r, 6) $row[5]
    Set-TextCell  $ws.Cells.Item($r, 7) $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: Update the "总计" sheet - insert a new top data row for 2022-Q4
# and shift the existing quarters down, renumbering the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$summary = @(
    @(0, "2022-Q4", 4, 1.03),
    @(1, "2022-Q3", 6, 0.9),
    @(2, "2022-Q1", 1, 0.05),
    @(3, "2021-Q2", 2, 0.68),
    @(4, "2021-Q1", 7, 0.84),
    @(5, "2020-Q4", 1, 0.31)
)

$r = 2
foreach ($row in $summary) {
    $idxCell = $total.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
